$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 196
$ws1.Range("F4").Value = 798
$ws1.Range("F6").Value = 14

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F4").Value = 196
$ws4.Range("F5").Value = 798
$ws4.Range("F7").Value = 14
